$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update the Role header cell (F13) and the per-account Role values (F14:F17) ---
# before deleting columns, so these become the merged "comma separated" role lists.
# (order of assignment controls the order new shared strings are appended in)
$ws.Range("F15").Value = "User, Administrator"
$ws.Range("F16").Value = "User, Administrator"
$ws.Range("F14").Value = "User, Administrator, ExcelImporter, ExecEngineer"
$ws.Range("F13").Value = "[Role,]"
$ws.Range("F17").Value = "User"

# --- 2. Move the cell comment from J12 to G12 (same text) before the columns shift ---
$oldComment = $ws.Range("J12").Comment
$commentText = $oldComment.Text()
$oldComment.Delete()

# --- 3. Delete the now redundant columns G:I (User/ExcelImporter/ExecEngineer columns) ---
# Column J (autoLoginAccount) and K (accIsGodAccount/misc) shift left into G and H.
$ws.Range("G1:I1").EntireColumn.Delete()

# --- 4. Re-create the comment at its new location (G12) ---
$newComment = $ws.Range("G12").AddComment($commentText)

# --- 5. Widen column F (G/H already inherit the old J/K widths automatically) ---
$ws.Columns.Item(6).ColumnWidth = 44.8

# --- 6. Update the selection shown in the sheet view ---
$ws.Range("I12:I13").Select()
